# RPA datasets push 2024-03-20
# Insert a new IPO record ("제일엠앤에스(구.제일기공)") as row 7 of Sheet1,
# pushing the existing rows 7-21 down to 8-22, and drop the oldest trailing
# record (the old row 21, "케이웨더") so the table stays a 20-row listing
# (A1:F21 overall incl. header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7..21 down by inserting a blank row at position 7.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row with the new offering's data.
$ws.Range("A7").Value = "제일엠앤에스(구.제일기공)"
$ws.Range("B7").Value = "2024.04.05~04.12"
$ws.Range("C7").Value = "15,000~18,000"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 36000
$ws.Range("F7").Value = "케이비증권"

# The insert pushed the old last row (케이웨더) down to row 22; remove it
# so the sheet keeps its original 21-row extent (A1:F21).
$ws.Rows("22:22").Delete()
